$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 8" shape (the TCP/connection-orientated summary box)
# by name rather than a hard-coded index, so the script is resilient to
# any shape re-ordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 8") {
        $shp = $candidate
    }
}

$tr = $shp.TextFrame.TextRange

# The author replaced "a host" with "a variety" in the middle of the
# paragraph (selected the words "a host " and typed "a variety " over
# them), which is what splits the single run into three runs.
$target = "a host "
$startPos = 1 + $tr.Text.IndexOf($target)
$selection = $tr.Characters($startPos, $target.Length)
$selection.Text = "a variety "

# The whole paragraph was also bumped from 48pt to 52pt.
$tr.Font.Size = 52

# With spAutoFit, PowerPoint relays out the text box to fit the new
# (larger) text; reproduce the resulting box height.
$shp.Height = 385.3265354330709
